# Insert a new data row at row 46 (pushing the existing rows 46-80 down to
# 47-81, same as Excel's native "Insert Sheet Rows" behaviour), then fill
# the newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(46).Insert()

$ws.Cells.Item(46, 1).Value  = 1
$ws.Cells.Item(46, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value  = 44574
$ws.Cells.Item(46, 5).Value  = 15
$ws.Cells.Item(46, 6).Value  = "Fruta"
$ws.Cells.Item(46, 7).Value  = 100102
$ws.Cells.Item(46, 8).Value  = "Cítricos"
$ws.Cells.Item(46, 9).Value  = 100102005
$ws.Cells.Item(46, 10).Value = "Naranja"
$ws.Cells.Item(46, 11).Value = "Navel Late"
$ws.Cells.Item(46, 12).Value = "Segunda"
$ws.Cells.Item(46, 13).Value = 250
$ws.Cells.Item(46, 14).Value = 800
$ws.Cells.Item(46, 15).Value = 850
$ws.Cells.Item(46, 16).Value = 825
$ws.Cells.Item(46, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(46, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(46, 19).Value = 825
$ws.Cells.Item(46, 20).Value = 1
